# Add team record (Wins/Losses/Ties) columns to the data sheet.
# Mirrors the author's commit: "Added team record to data" — W/L/T live on
# the same sheet (columns AD:AF) rather than a separate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 66

# New header cells, styled like the existing header row (bold, bordered,
# centered) by copying the format from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team record.
$ws.Range("AD2:AD" + $lastRow).Value = 67
$ws.Range("AE2:AE" + $lastRow).Value = 95
$ws.Range("AF2:AF" + $lastRow).Value = 0
